$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting rows 74:178 down to 75:179
$ws.Rows("74:74").Insert()

# Fill the constant columns (same values as every other data row in the block)
$ws.Range("A74").Value = 7
$ws.Range("B74").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C74").Value = "Ñuble"
$ws.Range("E74").Value = 16
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100103
$ws.Range("H74").Value = "Frutos de hueso (carozo)"
$ws.Range("I74").Value = 100103001
$ws.Range("J74").Value = "Cereza"
$ws.Range("T74").Value = 10

# Fill the new row's specific data
$ws.Range("D74").Value = 44524
$ws.Range("K74").Value = "Santina"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 60
$ws.Range("N74").Value = 21000
$ws.Range("O74").Value = 22000
$ws.Range("P74").Value = 21500
$ws.Range("Q74").Value = "$/bandeja 10 kilos"
$ws.Range("R74").Value = "Provincia de Curicó"
$ws.Range("S74").Value = 2150
